$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles + row height) of the last existing row (17) into
# the new row (18), then overwrite the values/text.
$ws.Range("A17:C17").Copy()
$ws.Range("A18:C18").PasteSpecial()

$ws.Range("A18").Value = "13. ja 14.8.2019"
$ws.Range("B18").Value = 12
$ws.Range("C18").Value = "Frontin ulkoasun parantelua, vanhaksi mennen Strava oauth tokenin uudistaminen"

$ws.Rows.Item(18).RowHeight = 26.65

# Scroll the view down one row and move the active selection to A19, as in
# the edited workbook.
$excel.ActiveWindow.ScrollRow = 14
$ws.Range("A19").Select()
